# Apply the "fixes in output dict creation, improved annotations, protocol
# changes" edit: rename the strength column header from "(raw)" to "(RMS)"
# and update the computed reaction/peak/difference/strength values for every
# experiment row (rows 16-18 were also reordered to reflect the corrected
# protocol grouping).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("E1").Value = "strength (RMS)"

# Row 2: April16_offset_PPI_18
$ws.Cells.Item(2,1).Value = "April16_offset_PPI_18"
$ws.Cells.Item(2,2).Value = 24.67
$ws.Cells.Item(2,3).Value = 30.67
$ws.Cells.Item(2,4).Value = 6
$ws.Cells.Item(2,5).Value = 103

# Row 3: April16_ASR_control
$ws.Cells.Item(3,1).Value = "April16_ASR_control"
$ws.Cells.Item(3,2).Value = 19
$ws.Cells.Item(3,3).Value = 27.5
$ws.Cells.Item(3,4).Value = 8.5
$ws.Cells.Item(3,5).Value = 69.25

# Row 4: April16_offset_PPI_12
$ws.Cells.Item(4,1).Value = "April16_offset_PPI_12"
$ws.Cells.Item(4,2).Value = 24.8
$ws.Cells.Item(4,3).Value = 30
$ws.Cells.Item(4,4).Value = 5.2
$ws.Cells.Item(4,5).Value = 82.40000000000001

# Row 5: April16_gap_depth
$ws.Cells.Item(5,1).Value = "April16_gap_depth"
$ws.Cells.Item(5,2).Value = 17.2
$ws.Cells.Item(5,3).Value = 30.4
$ws.Cells.Item(5,4).Value = 13.2
$ws.Cells.Item(5,5).Value = 108.4

# Row 6: April16_gap_duration_20
$ws.Cells.Item(6,1).Value = "April16_gap_duration_20"
$ws.Cells.Item(6,2).Value = 18
$ws.Cells.Item(6,3).Value = 30
$ws.Cells.Item(6,4).Value = 12
$ws.Cells.Item(6,5).Value = 96.75

# Row 7: April16_gap_duration_8
$ws.Cells.Item(7,1).Value = "April16_gap_duration_8"
$ws.Cells.Item(7,2).Value = 21.2
$ws.Cells.Item(7,3).Value = 26
$ws.Cells.Item(7,4).Value = 4.8
$ws.Cells.Item(7,5).Value = 97.40000000000001

# Row 8: April16_offset_PPI_20
$ws.Cells.Item(8,1).Value = "April16_offset_PPI_20"
$ws.Cells.Item(8,2).Value = 24
$ws.Cells.Item(8,3).Value = 30.8
$ws.Cells.Item(8,4).Value = 6.8
$ws.Cells.Item(8,5).Value = 52

# Row 9: April16_gap_duration_10
$ws.Cells.Item(9,1).Value = "April16_gap_duration_10"
$ws.Cells.Item(9,2).Value = 16.8
$ws.Cells.Item(9,3).Value = 28.4
$ws.Cells.Item(9,4).Value = 11.6
$ws.Cells.Item(9,5).Value = 75.8

# Row 10: April16_gap_duration_50
$ws.Cells.Item(10,1).Value = "April16_gap_duration_50"
$ws.Cells.Item(10,2).Value = 19.6
$ws.Cells.Item(10,3).Value = 30.4
$ws.Cells.Item(10,4).Value = 10.8
$ws.Cells.Item(10,5).Value = 90

# Row 11: April16_offset_PPI_50
$ws.Cells.Item(11,1).Value = "April16_offset_PPI_50"
$ws.Cells.Item(11,2).Value = 20.8
$ws.Cells.Item(11,3).Value = 30
$ws.Cells.Item(11,4).Value = 9.199999999999999
$ws.Cells.Item(11,5).Value = 79.2

# Row 12: April16_offset_PPI_14
$ws.Cells.Item(12,1).Value = "April16_offset_PPI_14"
$ws.Cells.Item(12,2).Value = 20.8
$ws.Cells.Item(12,3).Value = 30.8
$ws.Cells.Item(12,4).Value = 10
$ws.Cells.Item(12,5).Value = 86.2

# Row 13: April16_offset_PPI_16
$ws.Cells.Item(13,1).Value = "April16_offset_PPI_16"
$ws.Cells.Item(13,2).Value = 17
$ws.Cells.Item(13,3).Value = 30
$ws.Cells.Item(13,4).Value = 13
$ws.Cells.Item(13,5).Value = 104.25

# Row 14: April16_offset_PPI_6
$ws.Cells.Item(14,1).Value = "April16_offset_PPI_6"
$ws.Cells.Item(14,2).Value = 19
$ws.Cells.Item(14,3).Value = 30
$ws.Cells.Item(14,4).Value = 11
$ws.Cells.Item(14,5).Value = 93.25

# Row 15: April16_gap_duration_4
$ws.Cells.Item(15,1).Value = "April16_gap_duration_4"
$ws.Cells.Item(15,2).Value = 21.5
$ws.Cells.Item(15,3).Value = 33
$ws.Cells.Item(15,4).Value = 11.5
$ws.Cells.Item(15,5).Value = 93.5

# Row 16: April16_tone_in_noise (was April16_offset_PPI_10)
$ws.Cells.Item(16,1).Value = "April16_tone_in_noise"
$ws.Cells.Item(16,2).Value = 14.67
$ws.Cells.Item(16,3).Value = 30
$ws.Cells.Item(16,4).Value = 15.33
$ws.Cells.Item(16,5).Value = 107.33

# Row 17: April16_offset_PPI_10 (was April16_offset_PPI_4)
$ws.Cells.Item(17,1).Value = "April16_offset_PPI_10"
$ws.Cells.Item(17,2).Value = 30
$ws.Cells.Item(17,3).Value = 32
$ws.Cells.Item(17,4).Value = 2
$ws.Cells.Item(17,5).Value = 85

# Row 18: April16_offset_PPI_4 (was April16_tone_in_noise)
$ws.Cells.Item(18,1).Value = "April16_offset_PPI_4"
$ws.Cells.Item(18,2).Value = 23
$ws.Cells.Item(18,3).Value = 30.5
$ws.Cells.Item(18,4).Value = 7.5
$ws.Cells.Item(18,5).Value = 75.5

# Row 19: April16_offset_PPI_8
$ws.Cells.Item(19,1).Value = "April16_offset_PPI_8"
$ws.Cells.Item(19,2).Value = 18.5
$ws.Cells.Item(19,3).Value = 30
$ws.Cells.Item(19,4).Value = 11.5
$ws.Cells.Item(19,5).Value = 64.25
